# "Generate Report for Handback" - update localization-status report after
# the de-de handback file was regenerated and found in sync with en-US.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# --- Status / Error Detail / Handback DateTime updates ---------------------

# Status column: handback completed and is in sync with en-US source.
# (Overview!E2/F2 mirror the per-language Status cell for zh-cn/de-de.)
$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("C2").Value = "Handed back: in sync with en-US"
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"

# Error Detail column (P2): the "not latest" staleness warning is gone now
# that a fresh handback was generated.
$zhcn.Range("P2").Value = ""
$dede.Range("P2").Value = ""

# Latest Handback DateTime column (K2): stamp the new handback generation time.
$zhcn.Range("K2").Value = "2016-08-22 02:59:12"
$dede.Range("K2").Value = "2016-08-22 02:59:19"

# --- Column width adjustments (report columns widened / shrunk) ------------

$overview.Range("E1").ColumnWidth = 29.9777047293527
$overview.Range("F1").ColumnWidth = 29.9777047293527

$zhcn.Range("C1").ColumnWidth = 29.9777047293527
$zhcn.Range("P1").ColumnWidth = 13.7470528738839

$dede.Range("C1").ColumnWidth = 29.9777047293527
$dede.Range("P1").ColumnWidth = 13.7470528738839
